$d = $word.ActiveDocument

# --- Change 1: "Mode utilisation en énumération" paragraph ---
# Originally split into two runs separated by a _GoBack bookmark; the
# bookmark is removed and the text becomes one single run.
$t1 = "Mode utilisation en énumération => évolutif (ajout d’un mode admin par ex)"
$f1 = $d.Content.Find
$f1.Execute($t1, $true, $false, $false, $false, $false, $true, 1, $false, $t1, 2) | Out-Null

# --- Change 2: "CancelWorkTime" description paragraph ---
# Originally split into three runs around a proofErr gramStart/gramEnd
# pair wrapping the lone word "de"; becomes one single run.
$t2 = "Cette méthode permet de repasser en mode « consultation » et affecte le premier temps de travail de la tache courante au temps de travail courant."
$f2 = $d.Content.Find
$f2.Execute($t2, $true, $false, $false, $false, $false, $true, 1, $false, $t2, 2) | Out-Null

# --- Change 3: "On réinitialise ensuite la liste « " prefix ---
# Only the first three runs (before the spell-checked "ListTaskToAddOrDelete")
# are merged into one; the rest of the paragraph is untouched.
$t3 = "On réinitialise ensuite la liste « "
$f3 = $d.Content.Find
$f3.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3, 2) | Out-Null

# --- Change 4: trailing empty paragraphs + bookmark relocation ---
# Collapse the last two empty paragraphs into one, and put the _GoBack
# bookmark on the final (now last) paragraph of the document.
$count = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs($count - 1)
$delRange = $d.Range($secondToLast.Range.Start, $secondToLast.Range.End)
$delRange.Delete() | Out-Null

$newCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($newCount)
$d.Bookmarks.Add("_GoBack", $lastPara.Range) | Out-Null
